$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Cached "Date Placeholder" field text: "3/2/2019" -> "3/2/19".
#    This cached text lives on the Date Placeholder shape of the slide
#    master AND every slide layout (all of them share the same field id),
#    so touch each one that currently shows the old cached value.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "3/2/2019") {
                $tr.Text = "3/2/19"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 7, shape 8 ("矩形 10"): shrink the box and trim the leading
#    "本次作业内容较多，因些" clause from the first run, leaving the rest
#    of the paragraph ("copy & paste" / "一些基础组件" runs) untouched.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$sh8 = $s7.Shapes.Item(8)

$sh8.Width = 401.01834645669294

$tr8 = $sh8.TextFrame.TextRange
$firstRun = $tr8.Characters(1, 21)
$firstRun.Text = "建议可以从示例代码中"
